# Journal de travail - add new entries (week 4, 25-26 May 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (49) as the formatting template for the
# new rows, then overwrite the values so the existing cell styles
# (borders/number formats/alignment) are reused instead of new styles
# being created.
$template = $ws.Range("A49:F49")

$template.Copy()
$ws.Range("A50:F50").PasteSpecial(-4122)
$ws.Range("A50").Value = 45071
$ws.Range("B50").Value = 4
$ws.Range("C50").Value = 2.25
$ws.Range("D50").Value = "Documentation"
$ws.Range("E50").Value = "Avancer sur la documentation"

$template.Copy()
$ws.Range("A51:F51").PasteSpecial(-4122)
$ws.Range("A51").Value = 45071
$ws.Range("B51").Value = 4
$ws.Range("C51").Value = 1.5
$ws.Range("D51").Value = "Implémentation"
$ws.Range("E51").Value = "Corriger des erreurs d'affichage des informations de l'utilisateur"

$template.Copy()
$ws.Range("A52:F52").PasteSpecial(-4122)
$ws.Range("A52").Value = 45071
$ws.Range("B52").Value = 4
$ws.Range("C52").Value = 1.5
$ws.Range("D52").Value = "Implémentation"
$ws.Range("E52").Value = "Vérifier le fonctionnement générale du site "

$template.Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$ws.Range("A53").Value = 45072
$ws.Range("B53").Value = 4
$ws.Range("C53").Value = 1.5
$ws.Range("D53").Value = "Implémentation"
$ws.Range("E53").Value = "Ajouter du contenu à la page d'accueil et modifier la page boutique"

$template.Copy()
$ws.Range("A54:F54").PasteSpecial(-4122)
$ws.Range("A54").Value = 45072
$ws.Range("B54").Value = 4
$ws.Range("C54").Value = 2.25
$ws.Range("D54").Value = "Documentation"
$ws.Range("E54").Value = "Avancer sur la documentation"

$template.Copy()
$ws.Range("A55:F55").PasteSpecial(-4122)
$ws.Range("A55").Value = 45072
$ws.Range("B55").Value = 4
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = "Implémentation"
$ws.Range("E55").Value = "Mettre la nouvelle version sur SwissCenter"

# Rows 51 and 53 hold two-line wrapped descriptions, so Excel auto-grows
# them to 30pt (same as other wrapped rows already in the sheet).
$ws.Rows("51:51").RowHeight = 30
$ws.Rows("53:53").RowHeight = 30

# Match the final selection/scroll position left behind by the edit.
$ws.Range("F55").Select()
